$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Presentation / Pack Level'
$ws.Range("B1").Value = 'VMP / AMP'
$ws.Range("C1").Value = 'BNF Code'
$ws.Range("D1").Value = 'MDR: Product Description'
$ws.Range("E1").Value = 'SNOMED Code'
$ws.Range("F1").Value = 'DM+D: Product Description'
$ws.Range("G1").Value = 'DM+D:Product and Pack Description'

$ws.Range("A2").Value = 'Presentation'
$ws.Range("B2").Value = 'VMP'
$ws.Range("C2").Value = '''''0203020C0AAAAAA'
$ws.Range("D2").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E2").Value = '''''4744111000001109'
$ws.Range("F2").Value = 'Adenosine 6mg/2ml solution for injection vials'
$ws.Range("G2").Value = 'Adenosine 6mg/2ml solution for injection 6 vials'

$ws.Range("A3").Value = 'Pack'
$ws.Range("B3").Value = 'VMP'
$ws.Range("C3").Value = '''''0203020C0AAAAAA'
$ws.Range("D3").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E3").Value = '''''34516311000001106'
$ws.Range("F3").Value = 'Adenosine 6mg/2ml solution for injection vials'
$ws.Range("G3").Value = 'Adenosine 6mg/2ml solution for injection 5 vials'

$ws.Range("A4").Value = 'Presentation'
$ws.Range("B4").Value = 'AMP'
$ws.Range("C4").Value = '''''0203020C0BBAAAA'
$ws.Range("D4").Value = 'Adenocor_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E4").Value = '''''4744411000001104'
$ws.Range("F4").Value = 'Adenocor 6mg/2ml solution for injection vials (Sanofi)'

$ws.Range("A5").Value = 'Pack'
$ws.Range("B5").Value = 'AMP'
$ws.Range("C5").Value = '''''0203020C0BBAAAA'
$ws.Range("D5").Value = 'Adenocor_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E5").Value = '''''4744711000001105'
$ws.Range("F5").Value = 'Adenocor 6mg/2ml solution for injection vials (Sanofi)'
$ws.Range("G5").Value = 'Adenocor 6mg/2ml solution for injection (Sanofi) 6 vials'

$ws.Range("A6").Value = 'Presentation'
$ws.Range("B6").Value = 'AMP'
$ws.Range("C6").Value = '''''0203020C0AAAAAA'
$ws.Range("D6").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E6").Value = '''''19663311000001109'
$ws.Range("F6").Value = 'Adenosine 6mg/2ml solution for injection vials (Wockhardt UK Ltd)'

$ws.Range("A7").Value = 'Pack'
$ws.Range("B7").Value = 'AMP'
$ws.Range("C7").Value = '''''0203020C0AAAAAA'
$ws.Range("D7").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E7").Value = '''''19663411000001102'
$ws.Range("F7").Value = 'Adenosine 6mg/2ml solution for injection vials (Wockhardt UK Ltd)'
$ws.Range("G7").Value = 'Adenosine 6mg/2ml solution for injection (Wockhardt UK Ltd) 6 vials'

$ws.Range("A8").Value = 'Presentation'
$ws.Range("B8").Value = 'AMP'
$ws.Range("C8").Value = '''''0203020C0AAAAAA'
$ws.Range("D8").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E8").Value = '''''20009311000001102'
$ws.Range("F8").Value = 'Adenosine 6mg/2ml solution for injection vials (A A H Pharmaceuticals Ltd)'

$ws.Range("A9").Value = 'Pack'
$ws.Range("B9").Value = 'AMP'
$ws.Range("C9").Value = '''''0203020C0AAAAAA'
$ws.Range("D9").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E9").Value = '''''20009411000001109'
$ws.Range("F9").Value = 'Adenosine 6mg/2ml solution for injection vials (A A H Pharmaceuticals Ltd)'
$ws.Range("G9").Value = 'Adenosine 6mg/2ml solution for injection (A A H Pharmaceuticals Ltd) 6 vials'

$ws.Range("A10").Value = 'Presentation'
$ws.Range("B10").Value = 'AMP'
$ws.Range("C10").Value = '''''0203020C0AAAAAA'
$ws.Range("D10").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E10").Value = '''''21855411000001109'
$ws.Range("F10").Value = 'Adenosine 6mg/2ml solution for injection vials (Focus Pharmaceuticals Ltd)'

$ws.Range("A11").Value = 'Pack'
$ws.Range("B11").Value = 'AMP'
$ws.Range("C11").Value = '''''0203020C0AAAAAA'
$ws.Range("D11").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E11").Value = '''''21855511000001108'
$ws.Range("F11").Value = 'Adenosine 6mg/2ml solution for injection vials (Focus Pharmaceuticals Ltd)'
$ws.Range("G11").Value = 'Adenosine 6mg/2ml solution for injection (Focus Pharmaceuticals Ltd) 6 vials'

$ws.Range("A12").Value = 'Presentation'
$ws.Range("B12").Value = 'AMP'
$ws.Range("C12").Value = '''''0203020C0AAAAAA'
$ws.Range("D12").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E12").Value = '''''24530711000001102'
$ws.Range("F12").Value = 'Adenosine 6mg/2ml solution for injection vials (Alliance Healthcare (Distribution) Ltd)'

$ws.Range("A13").Value = 'Pack'
$ws.Range("B13").Value = 'AMP'
$ws.Range("C13").Value = '''''0203020C0AAAAAA'
$ws.Range("D13").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E13").Value = '''''24531011000001108'
$ws.Range("F13").Value = 'Adenosine 6mg/2ml solution for injection vials (Alliance Healthcare (Distribution) Ltd)'
$ws.Range("G13").Value = 'Adenosine 6mg/2ml solution for injection (Alliance Healthcare (Distribution) Ltd) 6 vials'

$ws.Range("A14").Value = 'Presentation'
$ws.Range("B14").Value = 'AMP'
$ws.Range("C14").Value = '''''0203020C0AAAAAA'
$ws.Range("D14").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E14").Value = '''''34516211000001103'
$ws.Range("F14").Value = 'Adenosine 6mg/2ml solution for injection vials (Peckforton Pharmaceuticals Ltd)'

$ws.Range("A15").Value = 'Pack'
$ws.Range("B15").Value = 'AMP'
$ws.Range("C15").Value = '''''0203020C0AAAAAA'
$ws.Range("D15").Value = 'Adenosine_I/V Inf 3mg/ml 2ml Vl'
$ws.Range("E15").Value = '''''34516411000001104'
$ws.Range("F15").Value = 'Adenosine 6mg/2ml solution for injection vials (Peckforton Pharmaceuticals Ltd)'
$ws.Range("G15").Value = 'Adenosine 6mg/2ml solution for injection (Peckforton Pharmaceuticals Ltd) 5 vials'

$ws.Range("A16").Value = 'Presentation'
$ws.Range("B16").Value = 'VMP'
$ws.Range("C16").Value = '''''1003020U0AAAIAI'
$ws.Range("D16").Value = 'Diclofenac Sod_Gel 2.32%'
$ws.Range("E16").Value = '''''22480211000001104'
$ws.Range("F16").Value = 'Diclofenac 2.32% gel'

$ws.Range("A17").Value = 'Pack'
$ws.Range("B17").Value = 'VMP'
$ws.Range("C17").Value = '''''1003020U0AAAIAI'
$ws.Range("D17").Value = 'Diclofenac Sod_Gel 2.32%'
$ws.Range("E17").Value = '''''22479411000001100'
$ws.Range("F17").Value = 'Diclofenac 2.32% gel'
$ws.Range("G17").Value = 'Diclofenac 2.32% 50 grams'

$ws.Range("A18").Value = 'Pack'
$ws.Range("B18").Value = 'VMP'
$ws.Range("C18").Value = '''''1003020U0AAAIAI'
$ws.Range("D18").Value = 'Diclofenac Sod_Gel 2.32%'
$ws.Range("E18").Value = '''''22479511000001101'
$ws.Range("F18").Value = 'Diclofenac 2.32% gel'
$ws.Range("G18").Value = 'Diclofenac 2.32% 30 grams'

$ws.Range("A19").Value = 'Pack'
$ws.Range("B19").Value = 'VMP'
$ws.Range("C19").Value = '''''1003020U0AAAIAI'
$ws.Range("D19").Value = 'Diclofenac Sod_Gel 2.32%'
$ws.Range("E19").Value = '''''26352411000001101'
$ws.Range("F19").Value = 'Diclofenac 2.32% gel'
$ws.Range("G19").Value = 'Diclofenac 2.32% 100 grams'

$ws.Range("A20").Value = 'Presentation'
$ws.Range("B20").Value = 'AMP'
$ws.Range("C20").Value = '''''1003020U0BBADAI'
$ws.Range("D20").Value = 'Voltarol 12 Hour Emulgel P_Gel 2.32%'
$ws.Range("E20").Value = '''''22479611000001102'
$ws.Range("F20").Value = 'Voltarol 12 Hour Emulgel P 2.32% gel (GlaxoSmithKline Consumer Healthcare)'

$ws.Range("A21").Value = 'Pack'
$ws.Range("B21").Value = 'AMP'
$ws.Range("C21").Value = '''''1003020U0BBADAI'
$ws.Range("D21").Value = 'Voltarol 12 Hour Emulgel P_Gel 2.32%'
$ws.Range("E21").Value = '''''22479711000001106'
$ws.Range("F21").Value = 'Voltarol 12 Hour Emulgel P 2.32% gel (GlaxoSmithKline Consumer Healthcare)'
$ws.Range("G21").Value = 'Voltarol 12 Hour Emulgel P 2.32% (GlaxoSmithKline Consumer Healthcare) 50 grams'

$ws.Range("A22").Value = 'Pack'
$ws.Range("B22").Value = 'AMP'
$ws.Range("C22").Value = '''''1003020U0BBADAI'
$ws.Range("D22").Value = 'Voltarol 12 Hour Emulgel P_Gel 2.32%'
$ws.Range("E22").Value = '''''22479911000001108'
$ws.Range("F22").Value = 'Voltarol 12 Hour Emulgel P 2.32% gel (GlaxoSmithKline Consumer Healthcare)'
$ws.Range("G22").Value = 'Voltarol 12 Hour Emulgel P 2.32% (GlaxoSmithKline Consumer Healthcare) 30 grams'

$ws.Range("A23").Value = 'Pack'
$ws.Range("B23").Value = 'AMP'
$ws.Range("C23").Value = '''''1003020U0BBADAI'
$ws.Range("D23").Value = 'Voltarol 12 Hour Emulgel P_Gel 2.32%'
$ws.Range("E23").Value = '''''26352611000001103'
$ws.Range("F23").Value = 'Voltarol 12 Hour Emulgel P 2.32% gel (GlaxoSmithKline Consumer Healthcare)'
$ws.Range("G23").Value = 'Voltarol 12 Hour Emulgel P 2.32% (GlaxoSmithKline Consumer Healthcare) 100 grams'

# Clear cells that no longer have values in the updated layout
$ws.Range("G4").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("G12").ClearContents()
$ws.Range("G14").ClearContents()

# Sheet dimension now spans to row 23

